# cryptos.xlsx - "Updated cryptos list" GitHub Actions refresh
# Rewrites the Coin / Link / Price / Volume(1h) figures for rows 2-51 of
# the active sheet to match the latest coinranking.com snapshot,
# including a rank swap between Monero and EthereumClassic (rows 26-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings such as "1.000" or "20.89" look like numbers to Excel's
# Value setter and would silently be parsed into doubles (dropping the
# trailing zero / changing the displayed text). Force those Price (column D)
# cells to Text first so they stay literal strings, matching how every
# other row in this column is already stored.
$textCells = @(
    'D5'
    'D7'
    'D8'
    'D9'
    'D10'
    'D12'
    'D13'
    'D14'
    'D16'
    'D17'
    'D18'
    'D19'
    'D20'
    'D21'
    'D24'
    'D25'
    'D26'
    'D27'
    'D28'
    'D30'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D36'
    'D37'
    'D38'
    'D39'
    'D40'
    'D41'
    'D43'
    'D44'
    'D45'
    'D46'
    'D47'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    'D2' = '27.558.46'
    'E2' = '  +2.29%  '
    'D3' = '1.790.09'
    'E3' = '  +4.09%  '
    'E4' = '  -0.05%  '
    'D5' = '313.89'
    'E5' = '  +1.36%  '
    'D7' = '0.5375'
    'E7' = '  +10.27%  '
    'D8' = '0.3765'
    'E8' = '  +7.88%  '
    'D9' = '42.92'
    'E9' = '  +2.03%  '
    'D10' = '0.07509'
    'E10' = '  +3.43%  '
    'E11' = '  +6.36%  '
    'D12' = '1.000'
    'E12' = '  -0.03%  '
    'D13' = '20.89'
    'E13' = '  +4.99%  '
    'D14' = '6.169'
    'E14' = '  +5.26%  '
    'D15' = '1.790.09'
    'E15' = '  +3.90%  '
    'D16' = '7.070'
    'E16' = '  +3.21%  '
    'D17' = '90.83'
    'E17' = '  +4.80%  '
    'D18' = '0.00001071'
    'E18' = '  +3.31%  '
    'D19' = '0.06498'
    'E19' = '  +1.94%  '
    'D20' = '0.9999'
    'E20' = '  -0.03%  '
    'D21' = '16.96'
    'E22' = '  +5.17%  '
    'D23' = '27.598.88'
    'E23' = '  +2.23%  '
    'D24' = '11.22'
    'E24' = '  +3.75%  '
    'D25' = '2.090'
    'E25' = '  +0.38%  '
    'B26' = 'Monero'
    'C26' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D26' = '155.57'
    'E26' = '  +1.03%  '
    'B27' = 'EthereumClassic'
    'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D27' = '20.47'
    'E27' = '  +2.78%  '
    'D28' = '2.379'
    'E28' = '  +15.15%  '
    'D29' = '1.995.40'
    'E29' = '  +4.10%  '
    'D30' = '121.89'
    'E30' = '  +0.96%  '
    'D31' = '1.123'
    'E31' = '  +9.40%  '
    'D32' = '0.1030'
    'E32' = '  +10.81%  '
    'D33' = '5.664'
    'E33' = '  +5.77%  '
    'D34' = '3.597'
    'E34' = '  +0.32%  '
    'D35' = '0.02290'
    'E35' = '  +5.23%  '
    'D36' = '8.722'
    'E36' = '  +16.14%  '
    'D37' = '0.06017'
    'E37' = '  +1.98%  '
    'D38' = '4.983'
    'E38' = '  +5.03%  '
    'D39' = '0.2088'
    'E39' = '  +4.65%  '
    'D40' = '11.41'
    'E40' = '  +3.93%  '
    'D41' = '0.6248'
    'E41' = '  +4.33%  '
    'E42' = '  -1.85%  '
    'D43' = '1.000'
    'E43' = '  +0.03%  '
    'D44' = '1.144'
    'E44' = '  +5.17%  '
    'D45' = '13.29'
    'E45' = '  +4.09%  '
    'D46' = '0.5863'
    'D47' = '3.632'
    'E47' = '  +1.58%  '
    'D48' = '121.60'
    'E48' = '  +3.24%  '
    'D49' = '1.913'
    'E49' = '  +4.27%  '
    'D50' = '1.134'
    'E50' = '  +2.12%  '
    'D51' = '0.06745'
    'E51' = '  +1.64%  '
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
